$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    # Set slide background to a solid white fill
    $bg = $s.Background
    $bg.Fill.Solid()
    $bg.Fill.ForeColor.RGB = 16777215

    # Style the title run: black text, Arial typeface
    $title = $s.Shapes.Item(1)
    $titleRange = $title.TextFrame.TextRange
    $titleRange.Font.Color.RGB = 0
    $titleRange.Font.Name = "Arial"
}
